$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case the " de "/"del"/"el"/"la"/"los"/"y" connector words in municipality (and state) names ---
$ws.Cells.Replace('Rincón de Romos', 'Rincón De Romos', 1) | Out-Null
$ws.Cells.Replace('San Francisco de los Romo', 'San Francisco De Los Romo', 1) | Out-Null
$ws.Cells.Replace('Playas de Rosarito', 'Playas De Rosarito', 1) | Out-Null
$ws.Cells.Replace('Comitán de Domínguez', 'Comitán De Domínguez', 1) | Out-Null
$ws.Cells.Replace('Mazapa de Madero', 'Mazapa De Madero', 1) | Out-Null
$ws.Cells.Replace('Hidalgo del Parral', 'Hidalgo Del Parral', 1) | Out-Null
$ws.Cells.Replace('Villa de Álvarez', 'Villa De Álvarez', 1) | Out-Null
$ws.Cells.Replace('Ciudad de México', 'Ciudad De México', 1) | Out-Null
$ws.Cells.Replace('Nombre de Dios', 'Nombre De Dios', 1) | Out-Null
$ws.Cells.Replace('Estado de México', 'Estado De México', 1) | Out-Null
$ws.Cells.Replace('Acambay de Ruíz Castañeda', 'Acambay De Ruíz Castañeda', 1) | Out-Null
$ws.Cells.Replace('Almoloya de Juárez', 'Almoloya De Juárez', 1) | Out-Null
$ws.Cells.Replace('Almoloya del Río', 'Almoloya Del Río', 1) | Out-Null
$ws.Cells.Replace('Atizapán de Zaragoza', 'Atizapán De Zaragoza', 1) | Out-Null
$ws.Cells.Replace('Chapa de Mota', 'Chapa De Mota', 1) | Out-Null
$ws.Cells.Replace('Ecatepec de Morelos', 'Ecatepec De Morelos', 1) | Out-Null
$ws.Cells.Replace('Ixtapan de la Sal', 'Ixtapan De La Sal', 1) | Out-Null
$ws.Cells.Replace('Naucalpan de Juárez', 'Naucalpan De Juárez', 1) | Out-Null
$ws.Cells.Replace('San Felipe del Progreso', 'San Felipe Del Progreso', 1) | Out-Null
$ws.Cells.Replace('Soyaniquilpan de Juárez', 'Soyaniquilpan De Juárez', 1) | Out-Null
$ws.Cells.Replace('Tenango del Valle', 'Tenango Del Valle', 1) | Out-Null
$ws.Cells.Replace('Tlalnepantla de Baz', 'Tlalnepantla De Baz', 1) | Out-Null
$ws.Cells.Replace('Villa del Carbón', 'Villa Del Carbón', 1) | Out-Null
$ws.Cells.Replace('San Miguel de Allende', 'San Miguel De Allende', 1) | Out-Null
$ws.Cells.Replace('Apaseo el Alto', 'Apaseo El Alto', 1) | Out-Null
$ws.Cells.Replace('Dolores Hidalgo Cuna de la Independencia Nacional', 'Dolores Hidalgo Cuna De La Independencia Nacional', 1) | Out-Null
$ws.Cells.Replace('Jaral del Progreso', 'Jaral Del Progreso', 1) | Out-Null
$ws.Cells.Replace('Purísima del Rincón', 'Purísima Del Rincón', 1) | Out-Null
$ws.Cells.Replace('San Diego de la Unión', 'San Diego De La Unión', 1) | Out-Null
$ws.Cells.Replace('San Francisco del Rincón', 'San Francisco Del Rincón', 1) | Out-Null
$ws.Cells.Replace('San Luis de la Paz', 'San Luis De La Paz', 1) | Out-Null
$ws.Cells.Replace('Santa Cruz de Juventino Rosas', 'Santa Cruz De Juventino Rosas', 1) | Out-Null
$ws.Cells.Replace('Silao de la Victoria', 'Silao De La Victoria', 1) | Out-Null
$ws.Cells.Replace('Valle de Santiago', 'Valle De Santiago', 1) | Out-Null
$ws.Cells.Replace('Acapulco de Juárez', 'Acapulco De Juárez', 1) | Out-Null
$ws.Cells.Replace('Ajuchitlán del Progreso', 'Ajuchitlán Del Progreso', 1) | Out-Null
$ws.Cells.Replace('Alcozauca de Guerrero', 'Alcozauca De Guerrero', 1) | Out-Null
$ws.Cells.Replace('Atoyac de Álvarez', 'Atoyac De Álvarez', 1) | Out-Null
$ws.Cells.Replace('Chilapa de Álvarez', 'Chilapa De Álvarez', 1) | Out-Null
$ws.Cells.Replace('Chilpancingo de los Bravo', 'Chilpancingo De Los Bravo', 1) | Out-Null
$ws.Cells.Replace('Coyuca de Benítez', 'Coyuca De Benítez', 1) | Out-Null
$ws.Cells.Replace('Coyuca de Catalán', 'Coyuca De Catalán', 1) | Out-Null
$ws.Cells.Replace('Huitzuco de los Figueroa', 'Huitzuco De Los Figueroa', 1) | Out-Null
$ws.Cells.Replace('Iguala de la Independencia', 'Iguala De La Independencia', 1) | Out-Null
$ws.Cells.Replace('Zihuatanejo de Azueta', 'Zihuatanejo De Azueta', 1) | Out-Null
$ws.Cells.Replace('La Unión de Isidoro Montes de Oca', 'La Unión De Isidoro Montes De Oca', 1) | Out-Null
$ws.Cells.Replace('Taxco de Alarcón', 'Taxco De Alarcón', 1) | Out-Null
$ws.Cells.Replace('Técpan de Galeana', 'Técpan De Galeana', 1) | Out-Null
$ws.Cells.Replace('Tepecoacuilco de Trujano', 'Tepecoacuilco De Trujano', 1) | Out-Null
$ws.Cells.Replace('Tlalixtaquilla de Maldonado', 'Tlalixtaquilla De Maldonado', 1) | Out-Null
$ws.Cells.Replace('Cuautepec de Hinojosa', 'Cuautepec De Hinojosa', 1) | Out-Null
$ws.Cells.Replace('Mineral de la Reforma', 'Mineral De La Reforma', 1) | Out-Null
$ws.Cells.Replace('Pachuca de Soto', 'Pachuca De Soto', 1) | Out-Null
$ws.Cells.Replace('Progreso de Obregón', 'Progreso De Obregón', 1) | Out-Null
$ws.Cells.Replace('Tepeji del Río de Ocampo', 'Tepeji Del Río De Ocampo', 1) | Out-Null
$ws.Cells.Replace('Tezontepec de Aldama', 'Tezontepec De Aldama', 1) | Out-Null
$ws.Cells.Replace('Tula de Allende', 'Tula De Allende', 1) | Out-Null
$ws.Cells.Replace('Zacualtipán de Ángeles', 'Zacualtipán De Ángeles', 1) | Out-Null
$ws.Cells.Replace('Acatlán de Juárez', 'Acatlán De Juárez', 1) | Out-Null
$ws.Cells.Replace('Atotonilco el Alto', 'Atotonilco El Alto', 1) | Out-Null
$ws.Cells.Replace('Autlán de Navarro', 'Autlán De Navarro', 1) | Out-Null
$ws.Cells.Replace('Concepción de Buenos Aires', 'Concepción De Buenos Aires', 1) | Out-Null
$ws.Cells.Replace('Cuautitlán de García Barragán', 'Cuautitlán De García Barragán', 1) | Out-Null
$ws.Cells.Replace('Huejuquilla el Alto', 'Huejuquilla El Alto', 1) | Out-Null
$ws.Cells.Replace('Lagos de Moreno', 'Lagos De Moreno', 1) | Out-Null
$ws.Cells.Replace('Ojuelos de Jalisco', 'Ojuelos De Jalisco', 1) | Out-Null
$ws.Cells.Replace('San Cristóbal de la Barranca', 'San Cristóbal De La Barranca', 1) | Out-Null
$ws.Cells.Replace('San Juan de los Lagos', 'San Juan De Los Lagos', 1) | Out-Null
$ws.Cells.Replace('San Juanito de Escobedo', 'San Juanito De Escobedo', 1) | Out-Null
$ws.Cells.Replace('San Martín de Bolaños', 'San Martín De Bolaños', 1) | Out-Null
$ws.Cells.Replace('San Miguel el Alto', 'San Miguel El Alto', 1) | Out-Null
$ws.Cells.Replace('San Sebastián del Oeste', 'San Sebastián Del Oeste', 1) | Out-Null
$ws.Cells.Replace('Santa María de los Ángeles', 'Santa María De Los Ángeles', 1) | Out-Null
$ws.Cells.Replace('Tamazula de Gordiano', 'Tamazula De Gordiano', 1) | Out-Null
$ws.Cells.Replace('Teocuitatlán de Corona', 'Teocuitatlán De Corona', 1) | Out-Null
$ws.Cells.Replace('Tepatitlán de Morelos', 'Tepatitlán De Morelos', 1) | Out-Null
$ws.Cells.Replace('Tizapán el Alto', 'Tizapán El Alto', 1) | Out-Null
$ws.Cells.Replace('Tlajomulco de Zúñiga', 'Tlajomulco De Zúñiga', 1) | Out-Null
$ws.Cells.Replace('Unión de Tula', 'Unión De Tula', 1) | Out-Null
$ws.Cells.Replace('Zacoalco de Torres', 'Zacoalco De Torres', 1) | Out-Null
$ws.Cells.Replace('Zapotlán el Grande', 'Zapotlán El Grande', 1) | Out-Null
$ws.Cells.Replace('Coalcomán de Vázquez Pallares', 'Coalcomán De Vázquez Pallares', 1) | Out-Null
$ws.Cells.Replace('Cojumatlán de Régules', 'Cojumatlán De Régules', 1) | Out-Null
$ws.Cells.Replace('Tiquicheo de Nicolás Romero', 'Tiquicheo De Nicolás Romero', 1) | Out-Null
$ws.Cells.Replace('Tetela del Volcán', 'Tetela Del Volcán', 1) | Out-Null
$ws.Cells.Replace('Amatlán de Cañas', 'Amatlán De Cañas', 1) | Out-Null
$ws.Cells.Replace('Bahía de Banderas', 'Bahía De Banderas', 1) | Out-Null
$ws.Cells.Replace('Ixtlán del Río', 'Ixtlán Del Río', 1) | Out-Null
$ws.Cells.Replace('Santa María del Oro', 'Santa María Del Oro', 1) | Out-Null
$ws.Cells.Replace('Ayoquezco de Aldama', 'Ayoquezco De Aldama', 1) | Out-Null
$ws.Cells.Replace('Heroica Ciudad de Huajuapan de León', 'Heroica Ciudad De Huajuapan De León', 1) | Out-Null
$ws.Cells.Replace('Heroica Ciudad de Tlaxiaco', 'Heroica Ciudad De Tlaxiaco', 1) | Out-Null
$ws.Cells.Replace('Heroica Ciudad de Juchitán de Zaragoza', 'Heroica Ciudad De Juchitán De Zaragoza', 1) | Out-Null
$ws.Cells.Replace('Miahuatlán de Porfirio Díaz', 'Miahuatlán De Porfirio Díaz', 1) | Out-Null
$ws.Cells.Replace('Oaxaca de Juárez', 'Oaxaca De Juárez', 1) | Out-Null
$ws.Cells.Replace('Ocotlán de Morelos', 'Ocotlán De Morelos', 1) | Out-Null
$ws.Cells.Replace('Putla Villa de Guerrero', 'Putla Villa De Guerrero', 1) | Out-Null
$ws.Cells.Replace('San Pedro el Alto', 'San Pedro El Alto', 1) | Out-Null
$ws.Cells.Replace('Teotitlán de Flores Magón', 'Teotitlán De Flores Magón', 1) | Out-Null
$ws.Cells.Replace('Tezoatlán de Segura y Luna', 'Tezoatlán De Segura Y Luna', 1) | Out-Null
$ws.Cells.Replace('Tlacolula de Matamoros', 'Tlacolula De Matamoros', 1) | Out-Null
$ws.Cells.Replace('Villa de Tututepec de Melchor Ocampo', 'Villa De Tututepec De Melchor Ocampo', 1) | Out-Null
$ws.Cells.Replace('Villa de Zaachila', 'Villa De Zaachila', 1) | Out-Null
$ws.Cells.Replace('Villa Sola de Vega', 'Villa Sola De Vega', 1) | Out-Null
$ws.Cells.Replace('Zapotitlán del Río', 'Zapotitlán Del Río', 1) | Out-Null
$ws.Cells.Replace('Los Reyes de Juárez', 'Los Reyes De Juárez', 1) | Out-Null
$ws.Cells.Replace('Palmar de Bravo', 'Palmar De Bravo', 1) | Out-Null
$ws.Cells.Replace('San Salvador el Seco', 'San Salvador El Seco', 1) | Out-Null
$ws.Cells.Replace('San Salvador el Verde', 'San Salvador El Verde', 1) | Out-Null
$ws.Cells.Replace('Tepanco de López', 'Tepanco De López', 1) | Out-Null
$ws.Cells.Replace('Cadereyta de Montes', 'Cadereyta De Montes', 1) | Out-Null
$ws.Cells.Replace('Jalpan de Serra', 'Jalpan De Serra', 1) | Out-Null
$ws.Cells.Replace('Landa de Matamoros', 'Landa De Matamoros', 1) | Out-Null
$ws.Cells.Replace('Pinal de Amoles', 'Pinal De Amoles', 1) | Out-Null
$ws.Cells.Replace('San Juan del Río', 'San Juan Del Río', 1) | Out-Null
$ws.Cells.Replace('Armadillo de los Infante', 'Armadillo De Los Infante', 1) | Out-Null
$ws.Cells.Replace('Mexquitic de Carmona', 'Mexquitic De Carmona', 1) | Out-Null
$ws.Cells.Replace('Santa María del Río', 'Santa María Del Río', 1) | Out-Null
$ws.Cells.Replace('Soledad de Graciano Sánchez', 'Soledad De Graciano Sánchez', 1) | Out-Null
$ws.Cells.Replace('Villa de Ramos', 'Villa De Ramos', 1) | Out-Null
$ws.Cells.Replace('Villa de Reyes', 'Villa De Reyes', 1) | Out-Null
$ws.Cells.Replace('Nacozari de García', 'Nacozari De García', 1) | Out-Null
$ws.Cells.Replace('Ixtacuixtla de Mariano Matamoros', 'Ixtacuixtla De Mariano Matamoros', 1) | Out-Null
$ws.Cells.Replace('Amatlán de los Reyes', 'Amatlán De Los Reyes', 1) | Out-Null
$ws.Cells.Replace('Castillo de Teayo', 'Castillo De Teayo', 1) | Out-Null
$ws.Cells.Replace('Cazones de Herrera', 'Cazones De Herrera', 1) | Out-Null
$ws.Cells.Replace('Poza Rica de Hidalgo', 'Poza Rica De Hidalgo', 1) | Out-Null
$ws.Cells.Replace('Concepción del Oro', 'Concepción Del Oro', 1) | Out-Null
$ws.Cells.Replace('Mezquital del Oro', 'Mezquital Del Oro', 1) | Out-Null
$ws.Cells.Replace('Moyahua de Estrada', 'Moyahua De Estrada', 1) | Out-Null
$ws.Cells.Replace('Nochistlán de Mejía', 'Nochistlán De Mejía', 1) | Out-Null
$ws.Cells.Replace('Teúl de González Ortega', 'Teúl De González Ortega', 1) | Out-Null
$ws.Cells.Replace('Tlaltenango de Sánchez Román', 'Tlaltenango De Sánchez Román', 1) | Out-Null
$ws.Cells.Replace('Villa de Cos', 'Villa De Cos', 1) | Out-Null

# --- Remove the trailing metadata/footer rows (713:717) and shrink the used range to A1:D711 ---
$ws.Rows("713:717").Delete()

Write-Host "Done"
